{"js": "const pairs = [\n  [\"2025-12-01 Monday\", \"2025-12-02 Tuesday\"],\n  [\"763\u00d77=5341\", \"202\u00d76=1212\"],\n  [\"544\u00d74=2176\", \"219\u00d76=1314\"],\n  [\"536\u00d77=3752\", \"877\u00d79=7893\"],\n  [\"369\u00d79=3321\", \"530\u00d79=4770\"],\n  [\"254\u00d75=1270\", \"798\u00d75=3990\"],\n  [\"325\u00d72=650\", \"722\u00d72=1444\"],\n  [\"253\u00d74=1012\", \"587\u00d72=1174\"],\n  [\"230\u00d73=690\", \"423\u00d76=2538\"],\n  [\"223\u00d78=1784\", \"224\u00d72=448\"],\n  [\"418\u00d79=3762\", \"619\u00d75=3095\"],\n  [\"207\u00d74=828\", \"931\u00d75=4655\"],\n  [\"988\u00d75=4940\", \"544\u00d79=4896\"],\n  [\"333\u00d72=666\", \"481\u00d79=4329\"],\n  [\"728\u00d75=3640\", \"625\u00d72=1250\"],\n  [\"843\u00d72=1686\", \"180\u00d72=360\"],\n  [\"463\u00d72=926\", \"458\u00d76=2748\"],\n  [\"746\u00d78=5968\", \"530\u00d79=4770\"],\n  [\"721\u00d78=5768\", \"309\u00d79=2781\"],\n  [\"582\u00d78=4656\", \"454\u00d75=2270\"],\n  [\"345\u00d79=3105\", \"130\u00d73=390\"],\n  [\"508\u00d75=2540\", \"494\u00d76=2964\"],\n  [\"616\u00d79=5544\", \"839\u00d73=2517\"],\n  [\"434\u00d77=3038\", \"354\u00d73=1062\"],\n  [\"308\u00d79=2772\", \"398\u00d73=1194\"],\n  [\"750\u00d74=3000\", \"673\u00d79=6057\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-12-01 Monday\", \"2025-12-02 Tuesday\")\n    ,@(\"763\u00d77=5341\", \"202\u00d76=1212\")\n    ,@(\"544\u00d74=2176\", \"219\u00d76=1314\")\n    ,@(\"536\u00d77=3752\", \"877\u00d79=7893\")\n    ,@(\"369\u00d79=3321\", \"530\u00d79=4770\")\n    ,@(\"254\u00d75=1270\", \"798\u00d75=3990\")\n    ,@(\"325\u00d72=650\", \"722\u00d72=1444\")\n    ,@(\"253\u00d74=1012\", \"587\u00d72=1174\")\n    ,@(\"230\u00d73=690\", \"423\u00d76=2538\")\n    ,@(\"223\u00d78=1784\", \"224\u00d72=448\")\n    ,@(\"418\u00d79=3762\", \"619\u00d75=3095\")\n    ,@(\"207\u00d74=828\", \"931\u00d75=4655\")\n    ,@(\"988\u00d75=4940\", \"544\u00d79=4896\")\n    ,@(\"333\u00d72=666\", \"481\u00d79=4329\")\n    ,@(\"728\u00d75=3640\", \"625\u00d72=1250\")\n    ,@(\"843\u00d72=1686\", \"180\u00d72=360\")\n    ,@(\"463\u00d72=926\", \"458\u00d76=2748\")\n    ,@(\"746\u00d78=5968\", \"530\u00d79=4770\")\n    ,@(\"721\u00d78=5768\", \"309\u00d79=2781\")\n    ,@(\"582\u00d78=4656\", \"454\u00d75=2270\")\n    ,@(\"345\u00d79=3105\", \"130\u00d73=390\")\n    ,@(\"508\u00d75=2540\", \"494\u00d76=2964\")\n    ,@(\"616\u00d79=5544\", \"839\u00d73=2517\")\n    ,@(\"434\u00d77=3038\", \"354\u00d73=1062\")\n    ,@(\"308\u00d79=2772\", \"398\u00d73=1194\")\n    ,@(\"750\u00d74=3000\", \"673\u00d79=6057\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
